$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.268.82'
$ws.Range("E2").Value = '  -2.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.566.81'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.61'
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  -5.07%  '
$ws.Range("E8").Value = '  -2.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0607'
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.83'
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.783.83'
$ws.Range("E12").Value = '  -3.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.572.83'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("E14").Value = '  -3.43%  '
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.275.44'
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '59.40'
$ws.Range("E17").Value = '  -2.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0714'
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '185.84'
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("E23").Value = '  -2.91%  '
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  -2.70%  '
$ws.Range("E27").Value = '  -7.21%  '
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("E29").Value = '  -3.69%  '
$ws.Range("E30").Value = '  -6.09%  '
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("E33").Value = '  -3.50%  '
$ws.Range("E34").Value = '  -1.57%  '
$ws.Range("E35").Value = '  -4.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.086.85'
$ws.Range("E36").Value = '  -3.52%  '
$ws.Range("E38").Value = '  -4.88%  '
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.773'
$ws.Range("E41").Value = '  -8.28%  '
$ws.Range("E42").Value = '  -0.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '93.33'
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.06'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.698.52'
$ws.Range("E45").Value = '  -3.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0107'
$ws.Range("E46").Value = '  -6.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.77'
$ws.Range("E47").Value = '  -3.10%  '
$ws.Range("E48").Value = '  -3.62%  '
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("E51").Value = '  -0.58%  '

Write-Host "Applied 69 cell updates"
